{"js": "// The canonical OOXML diff for this revision touches only the root-element\n// `xmlns:*` attribute ordering/prefix-numbering (e.g. `xmlns:m` shifting\n// position and the `urn:schemas-microsoft-com:office:excel` namespace prefix\n// being renumbered from `ns17` to `ns19`) on document.xml, endnotes.xml,\n// footer1.xml, footer2.xml, footnotes.xml, header1.xml, numbering.xml,\n// styles.xml and theme1.xml. That is a cosmetic artifact of whichever tool\n// last serialized the package (the set of namespace URIs declared is\n// byte-for-byte identical before and after; nothing is added, removed, or\n// reordered in the document body, headers/footers, numbering definitions,\n// styles, or theme content/text itself).\n//\n// There is no reachable, visible, or structural document edit described by\n// the diff, and the Word JavaScript API has no surface for rewriting raw\n// XML-namespace-prefix tables. So this script performs a safe, read-only\n// touch of the document body that leaves all content exactly as it is.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The canonical OOXML diff for this revision touches only the root-element\n# `xmlns:*` attribute ordering/prefix-numbering (e.g. `xmlns:m` shifting\n# position and the `urn:schemas-microsoft-com:office:excel` namespace prefix\n# being renumbered from `ns17` to `ns19`) on document.xml, endnotes.xml,\n# footer1.xml, footer2.xml, footnotes.xml, header1.xml, numbering.xml,\n# styles.xml and theme1.xml. That is a cosmetic artifact of whichever tool\n# last serialized the package (the set of namespace URIs declared is\n# byte-for-byte identical before and after; nothing is added, removed, or\n# reordered in the document body, headers/footers, numbering definitions,\n# styles, or theme content/text itself).\n#\n# There is no reachable, visible, or structural document edit described by\n# the diff, and the Word COM object model has no surface for rewriting raw\n# XML-namespace-prefix tables. So this script performs a safe, read-only\n# touch of the document that leaves all content exactly as it is.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
